# Update the answer table: replace each two-digit-division answer with the
# newly generated value. Cells are addressed by (row, column) in the single
# table on the page so that duplicate answer strings (e.g. "66÷5=13, 1")
# are not ambiguous.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "66÷5=13, 1"; New = "52÷3=17, 1" },
    @{ Row = 1;  Col = 2; Old = "90÷3=30, 0"; New = "58÷6=9, 4" },
    @{ Row = 1;  Col = 3; Old = "54÷7=7, 5";  New = "90÷6=15, 0" },
    @{ Row = 1;  Col = 4; Old = "87÷4=21, 3"; New = "23÷4=5, 3" },
    @{ Row = 1;  Col = 5; Old = "61÷5=12, 1"; New = "62÷9=6, 8" },

    @{ Row = 5;  Col = 1; Old = "46÷6=7, 4";  New = "83÷2=41, 1" },
    @{ Row = 5;  Col = 2; Old = "39÷6=6, 3";  New = "51÷8=6, 3" },
    @{ Row = 5;  Col = 4; Old = "38÷8=4, 6";  New = "66÷5=13, 1" },
    @{ Row = 5;  Col = 5; Old = "66÷2=33, 0"; New = "44÷5=8, 4" },

    @{ Row = 9;  Col = 1; Old = "50÷8=6, 2";  New = "48÷5=9, 3" },
    @{ Row = 9;  Col = 2; Old = "63÷9=7, 0";  New = "52÷3=17, 1" },
    @{ Row = 9;  Col = 3; Old = "53÷7=7, 4";  New = "18÷5=3, 3" },
    @{ Row = 9;  Col = 4; Old = "51÷6=8, 3";  New = "57÷7=8, 1" },
    @{ Row = 9;  Col = 5; Old = "40÷3=13, 1"; New = "98÷6=16, 2" },

    @{ Row = 13; Col = 1; Old = "10÷8=1, 2";  New = "53÷9=5, 8" },
    @{ Row = 13; Col = 2; Old = "56÷6=9, 2";  New = "91÷9=10, 1" },
    @{ Row = 13; Col = 3; Old = "76÷6=12, 4"; New = "23÷7=3, 2" },
    @{ Row = 13; Col = 4; Old = "64÷5=12, 4"; New = "88÷2=44, 0" },
    @{ Row = 13; Col = 5; Old = "54÷3=18, 0"; New = "34÷7=4, 6" },

    @{ Row = 17; Col = 1; Old = "67÷5=13, 2"; New = "52÷7=7, 3" },
    @{ Row = 17; Col = 2; Old = "49÷7=7, 0";  New = "10÷3=3, 1" },
    @{ Row = 17; Col = 3; Old = "92÷6=15, 2"; New = "99÷9=11, 0" },
    @{ Row = 17; Col = 4; Old = "93÷5=18, 3"; New = "10÷4=2, 2" },
    @{ Row = 17; Col = 5; Old = "13÷3=4, 1";  New = "71÷8=8, 7" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $range = $cell.Range
    # The cell's Range includes the trailing end-of-cell mark as a single
    # extra character position; moving the end back by one character unit
    # drops that mark so only the visible answer text remains selected.
    $range.MoveEnd(1, -1) | Out-Null
    if ($range.Text -eq $u.Old) {
        $range.Text = $u.New
    } else {
        $cell.Range.Find.Execute($u.Old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $u.New, 2) | Out-Null
    }
}
